$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns keep their original text (string) type
# by marking them as Text-formatted before assigning values, then resetting the style
# back to Normal so no visible/style difference is introduced versus the source file.
$numRange = $ws.Range("D2:E50")
$numRange.NumberFormat = "@"

$ws.Range("D2").Value = "69.020.14"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "3.803.33"
$ws.Range("E3").Value = "  +1.74%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "601.08"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("D6").Value = "163.32"
$ws.Range("E6").Value = "  -3.11%  "

$ws.Range("D7").Value = "3.799.83"
$ws.Range("E7").Value = "  +1.72%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -0.42%  "

$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("E11").Value = "  -1.82%  "

$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -0.30%  "

$ws.Range("D13").Value = "37.12"
$ws.Range("E13").Value = "  -2.96%  "

$ws.Range("D14").Value = "0.0000245"
$ws.Range("E14").Value = "  -1.36%  "

$ws.Range("D15").Value = "4.440.83"
$ws.Range("E15").Value = "  +1.76%  "

$ws.Range("D16").Value = "3.784.29"
$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("D17").Value = "69.148.04"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "7.43"
$ws.Range("E18").Value = "  +1.65%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.114"
$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "11.50"
$ws.Range("E20").Value = "  +5.78%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "17.20"
$ws.Range("E21").Value = "  +0.61%  "

$ws.Range("D22").Value = "485.89"
$ws.Range("E22").Value = "  -1.72%  "

$ws.Range("D23").Value = "0.719"
$ws.Range("E23").Value = "  -1.09%  "

$ws.Range("E24").Value = "  +4.51%  "

$ws.Range("D25").Value = "84.70"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").Value = "2.25"
$ws.Range("E26").Value = "  -3.41%  "

$ws.Range("D27").Value = "12.21"
$ws.Range("E27").Value = "  -1.22%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "10.00"
$ws.Range("E29").Value = "  -2.29%  "

$ws.Range("D30").Value = "2.97"
$ws.Range("E30").Value = "  -1.47%  "

$ws.Range("D31").Value = "7.99"
$ws.Range("E31").Value = "  -0.55%  "

$ws.Range("D32").Value = "2.38"
$ws.Range("E32").Value = "  -5.12%  "

$ws.Range("D33").Value = "3.954.95"
$ws.Range("E33").Value = "  +1.88%  "

$ws.Range("D34").Value = "31.65"
$ws.Range("E34").Value = "  -0.56%  "

$ws.Range("D35").Value = "3.747.58"
$ws.Range("E35").Value = "  +2.10%  "

$ws.Range("E36").Value = "  -1.86%  "

$ws.Range("D37").Value = "0.140"
$ws.Range("E37").Value = "  +4.95%  "

$ws.Range("E38").Value = "  +0.53%  "

$ws.Range("D39").Value = "5.88"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("E41").Value = "  +0.75%  "

$ws.Range("D42").Value = "0.319"
$ws.Range("E42").Value = "  -1.79%  "

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "48.55"
$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "426.03"
$ws.Range("E44").Value = "  -1.78%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "1.99"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").Value = "8.35"
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").Value = "2.824.02"
$ws.Range("E48").Value = "  +1.37%  "

$ws.Range("D49").Value = "141.47"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").Value = "39.31"
$ws.Range("E50").Value = "  -3.29%  "

# Reset style on the numeric-looking text columns back to Normal (removes the temporary
# Text number-format flag while keeping the values stored as text).
$numRange.Style = "Normal"
